$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style/format of row 314 down into the three new rows (315-317)
# so the new rows match the existing formatting (date style on column A, etc.)
$ws.Range("A314:G314").Copy($ws.Range("A315:G315"))
$ws.Range("A314:G314").Copy($ws.Range("A316:G316"))
$ws.Range("A314:G314").Copy($ws.Range("A317:G317"))

# Update existing row 314 (low/close changed)
$ws.Range("E314").Value = 106.43
$ws.Range("F314").Value = 109

# Row 315
$ws.Range("A315").Value = 45170.33333333334
$ws.Range("B315").Value = "FX_IDC:USDBDT"
$ws.Range("C315").Value = 109
$ws.Range("D315").Value = 110.28
$ws.Range("E315").Value = 108.53
$ws.Range("F315").Value = 109.97
$ws.Range("G315").Value = 0

# Row 316
$ws.Range("A316").Value = 45201.375
$ws.Range("B316").Value = "FX_IDC:USDBDT"
$ws.Range("C316").Value = 109.97
$ws.Range("D316").Value = 110.4
$ws.Range("E316").Value = 108.53
$ws.Range("F316").Value = 110
$ws.Range("G316").Value = 0

# Row 317
$ws.Range("A317").Value = 45231.375
$ws.Range("B317").Value = "FX_IDC:USDBDT"
$ws.Range("C317").Value = 110
$ws.Range("D317").Value = 110.63
$ws.Range("E317").Value = 108.5
$ws.Range("F317").Value = 110.48
$ws.Range("G317").Value = 0
